# Regenerate the "K" column (column G) values for rows 2-28 of the
# active worksheet. These values represent recalculated strikeout
# counts (K) replacing the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row (1-based Excel row) -> new K value
$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 2
    11 = 1
    12 = 2
    13 = 2
    14 = 5
    15 = 0
    16 = 3
    17 = 0
    18 = 2
    19 = 0
    20 = 2
    21 = 0
    22 = 2
    23 = 0
    24 = 0
    25 = 1
    26 = 4
    27 = 1
    28 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
